$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from an existing "Fecha" cell (D78) onto the
# still-empty D79:D92 range, then fill in the Avance (100) and Fecha values
# for rows 79-92, matching the progress made on video items 77-90.

$ws.Range("D78").Copy()
$ws.Range("D79:D92").PasteSpecial(-4122)

$ws.Range("C79:C92").Value = 100

$ws.Range("D79").Value = 44125
$ws.Range("D80").Value = 44125
$ws.Range("D81").Value = 44125
$ws.Range("D82").Value = 44126
$ws.Range("D83").Value = 44126
$ws.Range("D84").Value = 44126
$ws.Range("D85").Value = 44126
$ws.Range("D86").Value = 44127
$ws.Range("D87").Value = 44127
$ws.Range("D88").Value = 44127
$ws.Range("D89").Value = 44128
$ws.Range("D90").Value = 44128
$ws.Range("D91").Value = 44128
$ws.Range("D92").Value = 44129

# Carry over the stray formatted-but-empty cell style (matching the one
# already present at F6) to G92, left over from the user's editing/selection.
$ws.Range("F6").Copy()
$ws.Range("G92").PasteSpecial(-4122)

# Reflect where the user's cursor ended up after entering the last value.
$ws.Range("G92").Select()
